# "module 3 and 4"
#
# 1) Slide 7  (K-Nearest Neighbours (KNN)):
#       remove the empty, unused "Content Placeholder 5" shape.
# 2) Slide 12 (Support Vector Machines (SVM)):
#       remove the empty, unused "Content Placeholder 2" shape.
# 3) Slides 18, 19, 20 (Confusion Matrix: Another Example - three progressive
#    build copies of the same slide): fix a copy/paste mistake in the
#    hypothesis line, "Defendant" -> "Patient", so it reads
#    "H1 (Positive): Patient HAS a disease" (matching the H0 line right above
#    it, which already says "Patient").

$p = $ppt.ActivePresentation

# --- 1: drop the stray empty placeholder on slide 7 ------------------------
$slide7 = $p.Slides.Item(7)
for ($i = $slide7.Shapes.Count; $i -ge 1; $i--) {
    $shape = $slide7.Shapes.Item($i)
    if ($shape.Name -eq "Content Placeholder 5") {
        $shape.Delete()
    }
}

# --- 2: drop the stray empty placeholder on slide 12 ------------------------
$slide12 = $p.Slides.Item(12)
for ($i = $slide12.Shapes.Count; $i -ge 1; $i--) {
    $shape = $slide12.Shapes.Item($i)
    if ($shape.Name -eq "Content Placeholder 2") {
        $shape.Delete()
    }
}

# --- 3: "Defendant" -> "Patient" on the three confusion-matrix slides ------
$oldPhrase = "Defendant HAS a disease"
$oldWord = "Defendant"
$newWord = "Patient"

foreach ($slideIndex in 18, 19, 20) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            if ($para.Text.Contains($oldPhrase)) {
                # Replace just the word "Defendant" -> "Patient" in place,
                # using absolute character offsets into the shape's text
                # range so the untouched parts of the run are left alone.
                $wordStart = $para.Start + $para.Text.IndexOf($oldWord)
                $wordRange = $tr.Characters($wordStart, $oldWord.Length)
                $wordRange.Text = $newWord
            }
        }
    }
}

Write-Output "edit.ps1 applied"
